# Add the new "2022-Q1" sheet, positioned right after "2021-Q4" (i.e. right
# before "总计"), and populate it with the quarterly fund-holding detail
# rows. Then insert a new top data row into "总计" summarizing 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" worksheet right after "2021-Q4"
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add([Type]::Missing, $afterSheet)
$newSheet.Name = "2022-Q1"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("B1:H1").Style = "header"

# Data rows (A column holds a zero-based row index, B..H the fund detail)
$rows = @(
    @(0,  "006102", "浙商丰利增强债券",               "48.75", "47.92", "2.11", "1.0286", 7),
    @(1,  "688888", "浙商聚潮产业成长混合",           "8.25",  "93.40", "4.91", "0.4051", 8),
    @(2,  "010381", "浙商智选价值混合A",              "2.92",  "93.43", "4.63", "0.1352", 10),
    @(3,  "009246", "摩根士丹利华鑫ESG量化先行混合",  "4.04",  "92.17", "2.01", "0.0812", 4),
    @(4,  "460009", "华泰柏瑞量化先行混合A",          "9.13",  "90.47", "0.84", "0.0767", 10),
    @(5,  "011179", "浙商智选食品饮料股票A",          "0.22",  "91.35", "8.67", "0.0191", 1),
    @(6,  "012005", "信达澳银恒盛混合A",              "1.87",  "31.90", "0.92", "0.0172", 2),
    @(7,  "010382", "浙商智选价值混合C",              "0.34",  "93.43", "4.63", "0.0157", 10),
    @(8,  "009188", "鹏华股息精选混合",               "0.64",  "89.69", "1.86", "0.0119", 7),
    @(9,  "009658", "汇丰晋信中小盘低波动策略股票A",  "0.98",  "86.56", "1.08", "0.0106", 9),
    @(10, "013242", "北信瑞丰优势行业股票",           "0.82",  "92.63", "1.15", "0.0094", 8),
    @(11, "970073", "东证融汇成长优选混合A",          "0.68",  "82.02", "0.84", "0.0057", 10),
    @(12, "004352", "北信瑞丰研究精选股票",           "0.49",  "92.71", "1.10", "0.0054", 9),
    @(13, "011180", "浙商智选食品饮料股票C",          "0.05",  "91.35", "8.67", "0.0043", 1),
    @(14, "012006", "信达澳银恒盛混合C",              "0.31",  "31.90", "0.92", "0.0029", 2),
    @(15, "970074", "东证融汇成长优选混合C",          "0.27",  "82.02", "0.84", "0.0023", 10),
    @(16, "010246", "华泰柏瑞量化先行混合C",          "0.12",  "90.47", "0.84", "0.0010", 10),
    @(17, "009775", "汇丰晋信中小盘低波动策略股票C",  "0.04",  "86.56", "1.08", "0.0004", 9)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("A$r").Style = "header"
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("E$r").Value = $row[4]
    $newSheet.Range("F$r").Value = $row[5]
    $newSheet.Range("G$r").Value = $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------
# 2) Insert a new 2022-Q1 summary row at the top of "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("A2").EntireRow.Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A2").Style = "header"
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 18
$totalSheet.Range("D2").Value = 1.83

# Renumber the zero-based index column (A) for the rows that shifted down
for ($row = 3; $row -le 7; $row++) {
    $totalSheet.Cells.Item($row, 1).Value = $row - 2
}
